$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new habitat species data for rows 17 and 21 (species_codes column N)
$ws.Range("N17").Value = "EB, WCT"
$ws.Range("N21").Value = "EB, WCT"

# Add summary table data for hab_value (H) and upstream_habitat_length_m (M) on rows 23-24
$ws.Range("H23").Value = "high"
$ws.Range("M23").Value = 540
$ws.Range("H24").Value = "high"

# Update the view to match the saved selection/scroll position
$excel.Goto($ws.Range("M27"), $true)
$ws.Range("M27").Select()
